# Fix model emx files
# - Add a new "tags" worksheet (after "packages") describing the
#   biobankconnect tag references used by the "packages" sheet.
# - Update the "packages" sheet so its tags column (D2) references the new
#   tag identifiers instead of a raw hyperlinked URL.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "packages" sheet: replace the hyperlinked URL in D2 with
#    a comma separated list of tag identifiers, and drop the hyperlink.
# ---------------------------------------------------------------------
$packages = $wb.Worksheets("packages")
$packages.Range("D2").Hyperlinks.Delete()
$packages.Range("D2").Value = "miameenv_home,miameenv_pub1"
$packages.Range("D2").Font.Underline = $false
$packages.Range("D2").Font.Color = 0

# ---------------------------------------------------------------------
# 2. Add the new "tags" worksheet right after "packages" (i.e. as the
#    last sheet in the workbook).
# ---------------------------------------------------------------------
$tags = $wb.Worksheets.Add($null, $packages)
$tags.Name = "tags"

$tags.Range("A1").Value = "identifier"
$tags.Range("B1").Value = "objectIRI"
$tags.Range("C1").Value = "label"
$tags.Range("D1").Value = "relationLabel"
$tags.Range("E1").Value = "codeSystem"
$tags.Range("F1").Value = "relationIRI"

$tags.Range("A2").Value = "miameenv_home"
$tags.Range("B2").Value = "http://mibbi.sourceforge.net/projects/MIAME-Env.shtml"
$tags.Range("C2").Value = "http://mibbi.sourceforge.net/projects/MIAME-Env.shtml"
$tags.Range("D2").Value = "Homepage"
$tags.Range("E2").Value = "system"
$tags.Range("F2").Value = "http://molgenis.org/biobankconnect/link"

$tags.Range("A3").Value = "miameenv_pub1"
$tags.Range("B3").Value = "http://www.ncbi.nlm.nih.gov/pubmed/16901223"
$tags.Range("C3").Value = "http://www.ncbi.nlm.nih.gov/pubmed/16901223"
$tags.Range("D3").Value = "Publication"
$tags.Range("E3").Value = "system"
$tags.Range("F3").Value = "http://molgenis.org/biobankconnect/link"

$tags.Hyperlinks.Add($tags.Range("B2"), "http://mibbi.sourceforge.net/projects/MIAME-Env.shtml")
$tags.Hyperlinks.Add($tags.Range("C2"), "http://mibbi.sourceforge.net/projects/MIAME-Env.shtml")

$tags.Range("A1:F3").Select()

# ---------------------------------------------------------------------
# 3. Leave "packages" as the active sheet/selection, as in the source
#    workbook after the edit.
# ---------------------------------------------------------------------
$packages.Activate()
$packages.Range("D2").Select()
